$d = $word.ActiveDocument

$d.Content.Find.Execute("42-0=42", $true, $true, $false, $false, $false, $true, 1, $false, "44-29=15", 2) | Out-Null
$d.Content.Find.Execute("77-12=65", $true, $true, $false, $false, $false, $true, 1, $false, "47+22=69", 2) | Out-Null
$d.Content.Find.Execute("75-13=62", $true, $true, $false, $false, $false, $true, 1, $false, "93-78=15", 2) | Out-Null
$d.Content.Find.Execute("53+9=62", $true, $true, $false, $false, $false, $true, 1, $false, "66-56=10", 2) | Out-Null
$d.Content.Find.Execute("57-2=55", $true, $true, $false, $false, $false, $true, 1, $false, "54-18=36", 2) | Out-Null
$d.Content.Find.Execute("65-53=12", $true, $true, $false, $false, $false, $true, 1, $false, "80-55=25", 2) | Out-Null
$d.Content.Find.Execute("86-80=6", $true, $true, $false, $false, $false, $true, 1, $false, "42+39=81", 2) | Out-Null
$d.Content.Find.Execute("1+3=4", $true, $true, $false, $false, $false, $true, 1, $false, "73-40=33", 2) | Out-Null
$d.Content.Find.Execute("25+22=47", $true, $true, $false, $false, $false, $true, 1, $false, "57-24=33", 2) | Out-Null
$d.Content.Find.Execute("57+36=93", $true, $true, $false, $false, $false, $true, 1, $false, "84-54=30", 2) | Out-Null
$d.Content.Find.Execute("79-71=8", $true, $true, $false, $false, $false, $true, 1, $false, "45+49=94", 2) | Out-Null
$d.Content.Find.Execute("37+33=70", $true, $true, $false, $false, $false, $true, 1, $false, "49+7=56", 2) | Out-Null
$d.Content.Find.Execute("7+63=70", $true, $true, $false, $false, $false, $true, 1, $false, "82-53=29", 2) | Out-Null
$d.Content.Find.Execute("50-5=45", $true, $true, $false, $false, $false, $true, 1, $false, "52+11=63", 2) | Out-Null
$d.Content.Find.Execute("95-69=26", $true, $true, $false, $false, $false, $true, 1, $false, "16+39=55", 2) | Out-Null
$d.Content.Find.Execute("32+15=47", $true, $true, $false, $false, $false, $true, 1, $false, "81-72=9", 2) | Out-Null
$d.Content.Find.Execute("83-13=70", $true, $true, $false, $false, $false, $true, 1, $false, "3+63=66", 2) | Out-Null
$d.Content.Find.Execute("73-5=68", $true, $true, $false, $false, $false, $true, 1, $false, "45+16=61", 2) | Out-Null
$d.Content.Find.Execute("34+59=93", $true, $true, $false, $false, $false, $true, 1, $false, "30+35=65", 2) | Out-Null
$d.Content.Find.Execute("94-21=73", $true, $true, $false, $false, $false, $true, 1, $false, "65-42=23", 2) | Out-Null
$d.Content.Find.Execute("8+3=11", $true, $true, $false, $false, $false, $true, 1, $false, "17+18=35", 2) | Out-Null
$d.Content.Find.Execute("17+78=95", $true, $true, $false, $false, $false, $true, 1, $false, "15-11=4", 2) | Out-Null
$d.Content.Find.Execute("24-4=20", $true, $true, $false, $false, $false, $true, 1, $false, "34-7=27", 2) | Out-Null
$d.Content.Find.Execute("44+55=99", $true, $true, $false, $false, $false, $true, 1, $false, "88-60=28", 2) | Out-Null
$d.Content.Find.Execute("18-9=9", $true, $true, $false, $false, $false, $true, 1, $false, "12+46=58", 2) | Out-Null
$d.Content.Find.Execute("18+7=25", $true, $true, $false, $false, $false, $true, 1, $false, "99-89=10", 2) | Out-Null
$d.Content.Find.Execute("99-87=12", $true, $true, $false, $false, $false, $true, 1, $false, "40+2=42", 2) | Out-Null
$d.Content.Find.Execute("7+9=16", $true, $true, $false, $false, $false, $true, 1, $false, "34+31=65", 2) | Out-Null
$d.Content.Find.Execute("39-26=13", $true, $true, $false, $false, $false, $true, 1, $false, "94-81=13", 2) | Out-Null
$d.Content.Find.Execute("20-4=16", $true, $true, $false, $false, $false, $true, 1, $false, "58-48=10", 2) | Out-Null
$d.Content.Find.Execute("18+52=70", $true, $true, $false, $false, $false, $true, 1, $false, "24+24=48", 2) | Out-Null
$d.Content.Find.Execute("47-40=7", $true, $true, $false, $false, $false, $true, 1, $false, "9+58=67", 2) | Out-Null
$d.Content.Find.Execute("70+14=84", $true, $true, $false, $false, $false, $true, 1, $false, "56-3=53", 2) | Out-Null
$d.Content.Find.Execute("21+24=45", $true, $true, $false, $false, $false, $true, 1, $false, "6+56=62", 2) | Out-Null
$d.Content.Find.Execute("22+17=39", $true, $true, $false, $false, $false, $true, 1, $false, "88-67=21", 2) | Out-Null
$d.Content.Find.Execute("36+40=76", $true, $true, $false, $false, $false, $true, 1, $false, "81-19=62", 2) | Out-Null
$d.Content.Find.Execute("15+82=97", $true, $true, $false, $false, $false, $true, 1, $false, "91+6=97", 2) | Out-Null
$d.Content.Find.Execute("69-39=30", $true, $true, $false, $false, $false, $true, 1, $false, "93-70=23", 2) | Out-Null
$d.Content.Find.Execute("37-18=19", $true, $true, $false, $false, $false, $true, 1, $false, "84-51=33", 2) | Out-Null
$d.Content.Find.Execute("34+26=60", $true, $true, $false, $false, $false, $true, 1, $false, "86-5=81", 2) | Out-Null
$d.Content.Find.Execute("96+2=98", $true, $true, $false, $false, $false, $true, 1, $false, "49+21=70", 2) | Out-Null
$d.Content.Find.Execute("89-25=64", $true, $true, $false, $false, $false, $true, 1, $false, "99-83=16", 2) | Out-Null
$d.Content.Find.Execute("23+2=25", $true, $true, $false, $false, $false, $true, 1, $false, "96-17=79", 2) | Out-Null
$d.Content.Find.Execute("54-52=2", $true, $true, $false, $false, $false, $true, 1, $false, "58+6=64", 2) | Out-Null
$d.Content.Find.Execute("75-25=50", $true, $true, $false, $false, $false, $true, 1, $false, "54-4=50", 2) | Out-Null
$d.Content.Find.Execute("89-27=62", $true, $true, $false, $false, $false, $true, 1, $false, "96-29=67", 2) | Out-Null
$d.Content.Find.Execute("88-29=59", $true, $true, $false, $false, $false, $true, 1, $false, "82-67=15", 2) | Out-Null
$d.Content.Find.Execute("32-23=9", $true, $true, $false, $false, $false, $true, 1, $false, "48-32=16", 2) | Out-Null
$d.Content.Find.Execute("37-0=37", $true, $true, $false, $false, $false, $true, 1, $false, "86-83=3", 2) | Out-Null
$d.Content.Find.Execute("98-44=54", $true, $true, $false, $false, $false, $true, 1, $false, "81-3=78", 2) | Out-Null
$d.Content.Find.Execute("26+12=38", $true, $true, $false, $false, $false, $true, 1, $false, "86-68=18", 2) | Out-Null
$d.Content.Find.Execute("73+22=95", $true, $true, $false, $false, $false, $true, 1, $false, "56-12=44", 2) | Out-Null
$d.Content.Find.Execute("37+28=65", $true, $true, $false, $false, $false, $true, 1, $false, "63-27=36", 2) | Out-Null
$d.Content.Find.Execute("56+29=85", $true, $true, $false, $false, $false, $true, 1, $false, "61+38=99", 2) | Out-Null
$d.Content.Find.Execute("73-23=50", $true, $true, $false, $false, $false, $true, 1, $false, "5+19=24", 2) | Out-Null
$d.Content.Find.Execute("75-23=52", $true, $true, $false, $false, $false, $true, 1, $false, "61-55=6", 2) | Out-Null
$d.Content.Find.Execute("16-15=1", $true, $true, $false, $false, $false, $true, 1, $false, "54-31=23", 2) | Out-Null
$d.Content.Find.Execute("55-15=40", $true, $true, $false, $false, $false, $true, 1, $false, "15+14=29", 2) | Out-Null
$d.Content.Find.Execute("5+7=12", $true, $true, $false, $false, $false, $true, 1, $false, "42+3=45", 2) | Out-Null
$d.Content.Find.Execute("26+40=66", $true, $true, $false, $false, $false, $true, 1, $false, "3+23=26", 2) | Out-Null
$d.Content.Find.Execute("93-67=26", $true, $true, $false, $false, $false, $true, 1, $false, "84-36=48", 2) | Out-Null
$d.Content.Find.Execute("88-54=34", $true, $true, $false, $false, $false, $true, 1, $false, "48-6=42", 2) | Out-Null
$d.Content.Find.Execute("35+17=52", $true, $true, $false, $false, $false, $true, 1, $false, "32+12=44", 2) | Out-Null
$d.Content.Find.Execute("43-30=13", $true, $true, $false, $false, $false, $true, 1, $false, "98-39=59", 2) | Out-Null
$d.Content.Find.Execute("30-25=5", $true, $true, $false, $false, $false, $true, 1, $false, "42+44=86", 2) | Out-Null
$d.Content.Find.Execute("99-75=24", $true, $true, $false, $false, $false, $true, 1, $false, "28+40=68", 2) | Out-Null
$d.Content.Find.Execute("15-1=14", $true, $true, $false, $false, $false, $true, 1, $false, "0+98=98", 2) | Out-Null
$d.Content.Find.Execute("84-13=71", $true, $true, $false, $false, $false, $true, 1, $false, "73-52=21", 2) | Out-Null
$d.Content.Find.Execute("61-47=14", $true, $true, $false, $false, $false, $true, 1, $false, "90-11=79", 2) | Out-Null
$d.Content.Find.Execute("84-11=73", $true, $true, $false, $false, $false, $true, 1, $false, "49+13=62", 2) | Out-Null
$d.Content.Find.Execute("75+7=82", $true, $true, $false, $false, $false, $true, 1, $false, "93-56=37", 2) | Out-Null
$d.Content.Find.Execute("95-16=79", $true, $true, $false, $false, $false, $true, 1, $false, "25+8=33", 2) | Out-Null
$d.Content.Find.Execute("33-3=30", $true, $true, $false, $false, $false, $true, 1, $false, "51+4=55", 2) | Out-Null
$d.Content.Find.Execute("69-11=58", $true, $true, $false, $false, $false, $true, 1, $false, "18-1=17", 2) | Out-Null
$d.Content.Find.Execute("33+30=63", $true, $true, $false, $false, $false, $true, 1, $false, "24+57=81", 2) | Out-Null
$d.Content.Find.Execute("39-3=36", $true, $true, $false, $false, $false, $true, 1, $false, "66-61=5", 2) | Out-Null
$d.Content.Find.Execute("0+21=21", $true, $true, $false, $false, $false, $true, 1, $false, "8-2=6", 2) | Out-Null
$d.Content.Find.Execute("0+72=72", $true, $true, $false, $false, $false, $true, 1, $false, "56+30=86", 2) | Out-Null
$d.Content.Find.Execute("36+13=49", $true, $true, $false, $false, $false, $true, 1, $false, "92-81=11", 2) | Out-Null
$d.Content.Find.Execute("71-32=39", $true, $true, $false, $false, $false, $true, 1, $false, "47-3=44", 2) | Out-Null
$d.Content.Find.Execute("17-1=16", $true, $true, $false, $false, $false, $true, 1, $false, "80-45=35", 2) | Out-Null
$d.Content.Find.Execute("1+55=56", $true, $true, $false, $false, $false, $true, 1, $false, "6+43=49", 2) | Out-Null
$d.Content.Find.Execute("98-72=26", $true, $true, $false, $false, $false, $true, 1, $false, "54+7=61", 2) | Out-Null
$d.Content.Find.Execute("33+25=58", $true, $true, $false, $false, $false, $true, 1, $false, "89-87=2", 2) | Out-Null
$d.Content.Find.Execute("48+46=94", $true, $true, $false, $false, $false, $true, 1, $false, "39-38=1", 2) | Out-Null
$d.Content.Find.Execute("81-24=57", $true, $true, $false, $false, $false, $true, 1, $false, "58-36=22", 2) | Out-Null
$d.Content.Find.Execute("99-22=77", $true, $true, $false, $false, $false, $true, 1, $false, "29+3=32", 2) | Out-Null
$d.Content.Find.Execute("38-4=34", $true, $true, $false, $false, $false, $true, 1, $false, "54+22=76", 2) | Out-Null
$d.Content.Find.Execute("38-31=7", $true, $true, $false, $false, $false, $true, 1, $false, "72-59=13", 2) | Out-Null
$d.Content.Find.Execute("60+1=61", $true, $true, $false, $false, $false, $true, 1, $false, "79-46=33", 2) | Out-Null
$d.Content.Find.Execute("64-19=45", $true, $true, $false, $false, $false, $true, 1, $false, "40+26=66", 2) | Out-Null
$d.Content.Find.Execute("28+8=36", $true, $true, $false, $false, $false, $true, 1, $false, "64-41=23", 2) | Out-Null
$d.Content.Find.Execute("18+47=65", $true, $true, $false, $false, $false, $true, 1, $false, "84-51=33", 2) | Out-Null
$d.Content.Find.Execute("23-1=22", $true, $true, $false, $false, $false, $true, 1, $false, "32+62=94", 2) | Out-Null
$d.Content.Find.Execute("96-80=16", $true, $true, $false, $false, $false, $true, 1, $false, "31+56=87", 2) | Out-Null
$d.Content.Find.Execute("14+33=47", $true, $true, $false, $false, $false, $true, 1, $false, "34+37=71", 2) | Out-Null
$d.Content.Find.Execute("4+71=75", $true, $true, $false, $false, $false, $true, 1, $false, "44+26=70", 2) | Out-Null
$d.Content.Find.Execute("14+36=50", $true, $true, $false, $false, $false, $true, 1, $false, "25-15=10", 2) | Out-Null
$d.Content.Find.Execute("17+59=76", $true, $true, $false, $false, $false, $true, 1, $false, "3-0=3", 2) | Out-Null
$d.Content.Find.Execute("44+0=44", $true, $true, $false, $false, $false, $true, 1, $false, "1+87=88", 2) | Out-Null
